# "sixth commit selected travelers"
# - Re-book the return flight for a later date and give the two date columns
#   distinct "day month year" custom formats (double-digit day for the return
#   date, single-digit day for the departure date).
# - Drop the leftover explicit formatting that was sitting on the
#   BoardingPlace / LandingPlace / Adults / Child cells for row 2.
# - Leave the UI selection on H8 (last cell touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BoardingPlace (C2), LandingPlace (D2), Adults (G2), Child (H2): strip the
# stray style that was applied to these "selected travelers" cells.
$ws.Range("C2").ClearFormats()
$ws.Range("D2").ClearFormats()
$ws.Range("G2").ClearFormats()
$ws.Range("H2").ClearFormats()

# ReturnDate (F2): push the trip out to 12-Oct-2025 and format as "dd mmmm yyyy".
$ws.Range("F2").Value = 45942
$ws.Range("F2").NumberFormat = "[$-14009]dd\ mmmm\ yyyy;@"

# DepartureDate (E2): value is unchanged, but now rendered with a single-digit
# day ("d mmmm yyyy") instead of the zero-padded "dd mmmm yyyy".
$ws.Range("E2").NumberFormat = "[$-14009]d\ mmmm\ yyyy;@"

# Last UI action: select H8.
$ws.Range("H8").Select()
